# Crypto symbol list update — Sun Jan  8 13:36:06 UTC 2023 (GitHub Actions refresh)
# Updates Price (col D) and Volume(1h) (col E) text values for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $range = $ws.Range($cellRef)
    # Leading apostrophe forces text storage so numeric-looking/percent
    # strings (e.g. "261.84", "-2.22%") stay text, matching the source data.
    $range.Value = "'" + $newValue
    # Re-apply the default style so we do not leave a stray quote-prefix
    # number format behind on a cell that previously had none.
    $range.Style = "Normal"
}

Set-TextValue "D2" "261.84"
Set-TextValue "E3" "-2.22%"
Set-TextValue "D4" "4.698"
Set-TextValue "E4" "0.22%"
Set-TextValue "E5" "-0.74%"
Set-TextValue "D6" "6.702"
Set-TextValue "D7" "0.8513"
Set-TextValue "E7" "-0.36%"
Set-TextValue "D8" "0.9100"
Set-TextValue "E8" "-1.20%"
Set-TextValue "D9" "0.1403"
Set-TextValue "E9" "0.13%"
Set-TextValue "D10" "0.05104"
Set-TextValue "E10" "8.87%"
Set-TextValue "D11" "0.07098"
Set-TextValue "E11" "0.10%"
Set-TextValue "D12" "0.03114"
Set-TextValue "E12" "1.98%"
Set-TextValue "D13" "0.09035"
Set-TextValue "E13" "-0.28%"
Set-TextValue "D14" "0.001541"
Set-TextValue "E14" "0.75%"
Set-TextValue "D15" "0.0006155"
Set-TextValue "E15" "1.25%"
Set-TextValue "D16" "0.006008"
Set-TextValue "E16" "-0.61%"
Set-TextValue "D17" "3.450"
Set-TextValue "E17" "-0.04%"
Set-TextValue "D18" "3.167"
Set-TextValue "E18" "0.59%"
Set-TextValue "E19" "0.12%"
Set-TextValue "E21" "-0.35%"
Set-TextValue "D22" "4.112"
Set-TextValue "E22" "0.70%"
Set-TextValue "D23" "0.04243"
Set-TextValue "E23" "0.06%"
Set-TextValue "D24" "0.001179"
Set-TextValue "E24" "-3.04%"
Set-TextValue "E25" "6.92%"
Set-TextValue "E26" "0.05%"
Set-TextValue "E27" "23.08%"
Set-TextValue "D40" "0.03954"
Set-TextValue "E40" "2.23%"
Set-TextValue "D41" "0.1113"
Set-TextValue "E41" "-0.03%"
Set-TextValue "E42" "2.39%"
Set-TextValue "E43" "-14.78%"
Set-TextValue "E44" "-7.00%"
Set-TextValue "E45" "-0.61%"
Set-TextValue "E46" "0.07%"
Set-TextValue "D48" "0.2584"
Set-TextValue "E48" "62.26%"
Set-TextValue "E49" "0.07%"
Set-TextValue "E50" "0.07%"
